$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

# New row of results appended by the tracker's automatic update.
$ws.Cells.Item($row, 1).Value = 14637552
$ws.Cells.Item($row, 3).Value = "Henri Squire"
$ws.Cells.Item($row, 4).Value = "Niklas Schell"
$ws.Cells.Item($row, 5).Value = "Gana Niklas Schell"
$ws.Cells.Item($row, 6).Value = 6.5

# "fecha" is stored as plain text (e.g. "2025-09-08"), not a real date,
# matching the rest of the column - force text formatting so Excel
# doesn't auto-convert it to a date serial, then restore the default
# (unstyled) look used by every other row in the sheet.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2025-09-08"
$ws.Cells.Item($row, 2).Style = $ws.Cells.Item($row - 1, 2).Style
